# Generate Report for handoff
#
# - The previously "Ready for handoff" source file
#   8dab804a-4022-4673-99d5-1ef2a836b256.md was re-processed under a new
#   commit/guid (11be287a-2207-4f78-8175-ad03a7978651.md) with a refreshed
#   handoff package (new xlf hash, new handoff timestamps).
# - A brand new source file fbb329a8-c6cc-4ed8-b3f4-144a26eaf9ed.md showed up
#   whose handoff transform failed, so it is reported as "Ignored" on the
#   per-locale sheets (and "Handoff transform failed" on the Overview sheet),
#   inserted just above the always-last ".localization-config" row.

$wb = $excel.ActiveWorkbook

$oldGuid = "8dab804a-4022-4673-99d5-1ef2a836b256"
$newGuid = "11be287a-2207-4f78-8175-ad03a7978651"
$failedGuid = "fbb329a8-c6cc-4ed8-b3f4-144a26eaf9ed"
$newHash = "bbbd5e7cbba9e96f45b5c54c7db60319c975f813"

$newMdName = "$newGuid.md"
$failedMdName = "$failedGuid.md"
$cfgName = ".localization-config"

$newXlfZh = "$newGuid.$newHash.zh-cn.xlf"
$newXlfDe = "$newGuid.$newHash.de-de.xlf"

$newHandoffZh = "2016-01-17 10:30:24"
$newHandoffDe = "2016-01-17 10:30:35"
$epoch = "0001-01-01 00:00:00"

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/66b19eceff19f5170a95a1ce0c1c8fad4078650f/e2e/$newMdName"
$failedMdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/66b19eceff19f5170a95a1ce0c1c8fad4078650f/e2e/$failedMdName"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/66b19eceff19f5170a95a1ce0c1c8fad4078650f/.localization-config"
$xlfZhUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c66ceddc3ebc7940e0e55381620ea9e54eea2b1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/$newXlfZh"
$xlfDeUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fa97e90e15d490e2c598a3f88e5384c65aef5bd7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/$newXlfDe"

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("B2").Value = "Ready for handoff"
$wsOverview.Range("C2").Value = "Ready for handoff"

$wsOverview.Range("A3").Value = $failedMdName
$wsOverview.Range("B3").Value = "Handoff transform failed"
$wsOverview.Range("C3").Value = "Handoff transform failed"

$wsOverview.Range("A4").Value = $cfgName
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdUrl, "", "", $newMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $failedMdUrl, "", "", $failedMdName)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("B2").Value = "Ready for handoff"
$wsZh.Range("C2").Value = $newXlfZh
$wsZh.Range("D2").Value = $newHandoffZh
$wsZh.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G2").Value = $epoch
$wsZh.Range("H2").Value = "Include"

$wsZh.Range("A3").Value = $failedMdName
$wsZh.Range("B3").Value = "Handoff transform failed"
$wsZh.Range("D3").Value = $epoch
$wsZh.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G3").Value = $epoch
$wsZh.Range("H3").Value = "Ignored"

$wsZh.Range("A4").Value = $cfgName
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = $epoch
$wsZh.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("G4").Value = $epoch
$wsZh.Range("H4").Value = "Ignored"

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl, "", "", $newMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $xlfZhUrl, "", "", $newXlfZh)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $failedMdUrl, "", "", $failedMdName)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $cfgUrl, "", "", $cfgName)

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("B2").Value = "Ready for handoff"
$wsDe.Range("C2").Value = $newXlfDe
$wsDe.Range("D2").Value = $newHandoffDe
$wsDe.Range("D2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G2").Value = $epoch
$wsDe.Range("H2").Value = "Include"

$wsDe.Range("A3").Value = $failedMdName
$wsDe.Range("B3").Value = "Handoff transform failed"
$wsDe.Range("D3").Value = $epoch
$wsDe.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G3").Value = $epoch
$wsDe.Range("H3").Value = "Ignored"

$wsDe.Range("A4").Value = $cfgName
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = $epoch
$wsDe.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("G4").Value = $epoch
$wsDe.Range("H4").Value = "Ignored"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl, "", "", $newMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $xlfDeUrl, "", "", $newXlfDe)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $failedMdUrl, "", "", $failedMdName)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $cfgUrl, "", "", $cfgName)

Write-Host "Report regenerated for handoff"
